$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 8087.0586
$ws.Range("I11").Value = 8087.0586
$ws.Range("K11").Value = 8087.0586
$ws.Range("M11").Value = -7947.0586

$ws.Range("H18").Value = 539
$ws.Range("I18").Value = 539
$ws.Range("K18").Value = 539
$ws.Range("M18").Value = -255

$ws.Range("H28").Value = 243.71428
$ws.Range("I28").Value = 184.33333
$ws.Range("K28").Value = 184.33333
$ws.Range("M28").Value = 300.66667

$ws.Range("H86").Value = 35716856
$ws.Range("I86").Value = 83335416
$ws.Range("J86").Value = 2934.8125
$ws.Range("K86").Value = 83335416
$ws.Range("L86").Value = 2934.8125
$ws.Range("M86").Value = -83334293
$ws.Range("N86").Value = -5180.8125

$ws.Range("H89").Value = 35716856
$ws.Range("I89").Value = 83335416
$ws.Range("J89").Value = 2934.8125
$ws.Range("K89").Value = 416677080
$ws.Range("L89").Value = 14674.0625
$ws.Range("M89").Value = -416671464
$ws.Range("N89").Value = -25906.0625

$ws.Range("H92").Value = 1512908.1
$ws.Range("I92").Value = 625806.5
$ws.Range("J92").Value = 5209165
$ws.Range("K92").Value = 625806.5
$ws.Range("L92").Value = 5209165
$ws.Range("M92").Value = -624558.5
$ws.Range("N92").Value = -5211661

$ws.Range("H132").Value = 4465596.5
$ws.Range("I132").Value = 5556626
$ws.Range("K132").Value = 16669878
$ws.Range("M132").Value = -16667348

$ws.Range("H137").Value = 2089471
$ws.Range("I137").Value = 6035.96
$ws.Range("J137").Value = 4354074.5
$ws.Range("K137").Value = 18107.88
$ws.Range("L137").Value = 13062223.5
$ws.Range("M137").Value = -15557.88
$ws.Range("N137").Value = -13067323.5

$ws.Range("H138").Value = 5138.8125
$ws.Range("I138").Value = 5893.35
$ws.Range("K138").Value = 17680.05
$ws.Range("M138").Value = -12540.05

$ws.Range("H141").Value = 10527.767
$ws.Range("I141").Value = 2675.7896
$ws.Range("J141").Value = 24090.273
$ws.Range("K141").Value = 8027.3688
$ws.Range("L141").Value = 72270.819
$ws.Range("M141").Value = -2847.3688
$ws.Range("N141").Value = -82630.819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 55000
$ws.Range("J44").Value = 55000
$ws.Range("L44").Value = 55000
$ws.Range("N44").Value = -55976

$ws.Range("H61").Value = 1463592.9
$ws.Range("I61").Value = 54710.65
$ws.Range("K61").Value = 54710.65
$ws.Range("M61").Value = -54498.65

$ws.Range("H74").Value = 287517.8
$ws.Range("I74").Value = 1324.2373
$ws.Range("J74").Value = 1176224.1
$ws.Range("K74").Value = 1324.2373
$ws.Range("L74").Value = 1176224.1
$ws.Range("M74").Value = -450.2373
$ws.Range("N74").Value = -1177972.1

$ws.Range("H77").Value = 287517.8
$ws.Range("I77").Value = 1324.2373
$ws.Range("J77").Value = 1176224.1
$ws.Range("K77").Value = 6621.1865
$ws.Range("L77").Value = 5881120.5
$ws.Range("M77").Value = -2253.1865
$ws.Range("N77").Value = -5889856.5

$ws.Range("H92").Value = 57800
$ws.Range("J92").Value = 57800
$ws.Range("L92").Value = 57800
$ws.Range("N92").Value = -62792

$ws.Range("H122").Value = 2751.0312
$ws.Range("I122").Value = 2345.2917
$ws.Range("J122").Value = 3968.25
$ws.Range("K122").Value = 7035.875100000001
$ws.Range("L122").Value = 11904.75
$ws.Range("M122").Value = -4585.875100000001
$ws.Range("N122").Value = -16804.75

$ws.Range("H132").Value = 4262.45
$ws.Range("I132").Value = 3899.2856
$ws.Range("J132").Value = 5109.8335
$ws.Range("K132").Value = 11697.8568
$ws.Range("L132").Value = 15329.5005
$ws.Range("M132").Value = -9167.856800000001
$ws.Range("N132").Value = -20389.5005

$ws.Range("H136").Value = 1463592.9
$ws.Range("I136").Value = 54710.65
$ws.Range("K136").Value = 164131.95
$ws.Range("M136").Value = -161581.95

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 500
$ws.Range("I11").Value = 500
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 500
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -360
$ws.Range("N11").ClearContents()

$ws.Range("H22").Value = 1224.6
$ws.Range("I22").Value = 1010
$ws.Range("J22").Value = 1725.3334
$ws.Range("K22").Value = 1010
$ws.Range("L22").Value = 1725.3334
$ws.Range("M22").Value = -837
$ws.Range("N22").Value = -2071.3334

$ws.Range("H80").Value = 45455124
$ws.Range("J80").Value = 358.92856
$ws.Range("L80").Value = 358.92856
$ws.Range("N80").Value = -2354.92856

$ws.Range("H83").Value = 45455124
$ws.Range("J83").Value = 358.92856
$ws.Range("L83").Value = 1794.6428
$ws.Range("N83").Value = -11778.6428

$ws.Range("H99").Value = 11859.218
$ws.Range("I99").Value = 13805.934
$ws.Range("K99").Value = 13805.934
$ws.Range("M99").Value = -12307.934

$ws.Range("H105").Value = 4485.1133
$ws.Range("I105").Value = 4580.3076
$ws.Range("K105").Value = 4580.3076
$ws.Range("M105").Value = -2833.3076

$ws.Range("H134").Value = 19151086
$ws.Range("I134").Value = 1695.6666
$ws.Range("K134").Value = 5086.9998
$ws.Range("M134").Value = -2551.9998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 1880047.2
$ws.Range("I3").Value = 2506563.2
$ws.Range("K3").Value = 2506563.2
$ws.Range("M3").Value = -2506450.2

$ws.Range("H16").Value = 8941018
$ws.Range("I16").Value = 17863538
$ws.Range("K16").Value = 17863538
$ws.Range("M16").Value = -17863251

$ws.Range("H31").Value = 2112.8135
$ws.Range("I31").Value = 2998.889
$ws.Range("J31").Value = 1953.32
$ws.Range("K31").Value = 2998.889
$ws.Range("L31").Value = 1953.32
$ws.Range("M31").Value = -2703.889
$ws.Range("N31").Value = -2543.32

$ws.Range("H34").Value = 2112.8135
$ws.Range("I34").Value = 2998.889
$ws.Range("J34").Value = 1953.32
$ws.Range("K34").Value = 2998.889
$ws.Range("L34").Value = 1953.32
$ws.Range("M34").Value = -2796.889
$ws.Range("N34").Value = -2357.32

$ws.Range("H58").Value = 5136.926
$ws.Range("I58").Value = 9935.546
$ws.Range("J58").Value = 1837.875
$ws.Range("K58").Value = 9935.546
$ws.Range("L58").Value = 1837.875
$ws.Range("M58").Value = -9732.546
$ws.Range("N58").Value = -2243.875

$ws.Range("H62").Value = 5862.9
$ws.Range("J62").Value = 3943.5
$ws.Range("L62").Value = 3943.5
$ws.Range("N62").Value = -5191.5

$ws.Range("H65").Value = 5862.9
$ws.Range("J65").Value = 3943.5
$ws.Range("L65").Value = 19717.5
$ws.Range("N65").Value = -25957.5

$ws.Range("H107").Value = 1512.8889
$ws.Range("I107").Value = 1280.2
$ws.Range("J107").Value = 1803.75
$ws.Range("K107").Value = 1280.2
$ws.Range("L107").Value = 1803.75
$ws.Range("M107").Value = 639.8
$ws.Range("N107").Value = -5643.75

$ws.Range("H113").Value = 8941018
$ws.Range("I113").Value = 17863538
$ws.Range("K113").Value = 17863538
$ws.Range("M113").Value = -17861368

$ws.Range("H132").Value = 19610908
$ws.Range("I132").Value = 2669.7273
$ws.Range("J132").Value = 55559344
$ws.Range("K132").Value = 8009.1819
$ws.Range("L132").Value = 166678032
$ws.Range("M132").Value = -5479.1819
$ws.Range("N132").Value = -166683092

$ws.Range("H134").Value = 2275.88
$ws.Range("I134").Value = 2036.9375
$ws.Range("K134").Value = 6110.8125
$ws.Range("M134").Value = -3575.8125

$ws.Range("H136").Value = 5136.926
$ws.Range("I136").Value = 9935.546
$ws.Range("J136").Value = 1837.875
$ws.Range("K136").Value = 29806.638
$ws.Range("L136").Value = 5513.625
$ws.Range("M136").Value = -27256.638
$ws.Range("N136").Value = -10613.625

$ws.Range("H141").Value = 211611.1
$ws.Range("I141").Value = 72115
$ws.Range("K141").Value = 72115
$ws.Range("M141").Value = -66935

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 3168.625
$ws.Range("I64").Value = 3042.625
$ws.Range("J64").Value = 3294.625
$ws.Range("K64").Value = 9127.875
$ws.Range("L64").Value = 9883.875
$ws.Range("M64").Value = -8857.875
$ws.Range("N64").Value = -10423.875

$ws.Range("H67").Value = 3168.625
$ws.Range("I67").Value = 3042.625
$ws.Range("J67").Value = 3294.625
$ws.Range("K67").Value = 9127.875
$ws.Range("L67").Value = 9883.875
$ws.Range("M67").Value = -8191.875
$ws.Range("N67").Value = -11755.875

$ws.Range("H92").Value = 520.5714
$ws.Range("I92").Value = 537
$ws.Range("J92").Value = 498.66666
$ws.Range("K92").Value = 1611
$ws.Range("L92").Value = 1495.99998
$ws.Range("M92").Value = -363
$ws.Range("N92").Value = -3991.99998

$ws.Range("H133").Value = 4340.6665
$ws.Range("J133").Value = 11500
$ws.Range("L133").Value = 34500
$ws.Range("N133").Value = -44620

$ws.Range("H138").Value = 5050.2666
$ws.Range("I138").Value = 5187.9165
$ws.Range("J138").Value = 4499.6665
$ws.Range("K138").Value = 15563.7495
$ws.Range("L138").Value = 13498.9995
$ws.Range("M138").Value = -10423.7495
$ws.Range("N138").Value = -23778.9995

$ws.Range("H140").Value = 1445.4166
$ws.Range("I140").Value = 891.75
$ws.Range("J140").Value = 5874.75
$ws.Range("K140").Value = 2675.25
$ws.Range("L140").Value = 17624.25
$ws.Range("M140").Value = 2504.75
$ws.Range("N140").Value = -27984.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 94.09524
$ws.Range("I2").Value = 66.73333
$ws.Range("J2").Value = 162.5
$ws.Range("K2").Value = 66.73333
$ws.Range("L2").Value = 162.5
$ws.Range("M2").Value = 46.26667
$ws.Range("N2").Value = -388.5

$ws.Range("H4").Value = 200
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

$ws.Range("H68").Value = 55555
$ws.Range("J68").Value = 55555
$ws.Range("L68").Value = 55555
$ws.Range("N68").Value = -57177

$ws.Range("H70").Value = 5450
$ws.Range("I70").Value = 4000
$ws.Range("J70").Value = 6900
$ws.Range("K70").Value = 4000
$ws.Range("L70").Value = 6900
$ws.Range("N70").Value = -7440
$ws.Range("M70").Value = -3730

$ws.Range("H71").Value = 55555
$ws.Range("J71").Value = 55555
$ws.Range("L71").Value = 166665
$ws.Range("N71").Value = -174777

$ws.Range("H73").Value = 5450
$ws.Range("I73").Value = 4000
$ws.Range("J73").Value = 6900
$ws.Range("K73").Value = 4000
$ws.Range("L73").Value = 6900
$ws.Range("N73").Value = -8772
$ws.Range("M73").Value = -3064

$ws.Range("H74").Value = 49999.5
$ws.Range("J74").Value = 49999.5
$ws.Range("L74").Value = 49999.5
$ws.Range("N74").Value = -51871.5

$ws.Range("H77").Value = 49999.5
$ws.Range("J77").Value = 49999.5
$ws.Range("L77").Value = 149998.5
$ws.Range("N77").Value = -159358.5

$ws.Range("H97").Value = 1226.75
$ws.Range("I97").Value = 1226.75
$ws.Range("K97").Value = 1226.75
$ws.Range("M97").Value = -730.75

$ws.Range("H102").Value = 71430930
$ws.Range("I102").Value = 83335250
$ws.Range("K102").Value = 83335250
$ws.Range("M102").Value = -83333628

$ws.Range("H113").Value = 2134.4443
$ws.Range("J113").Value = 1900
$ws.Range("L113").Value = 1900
$ws.Range("N113").Value = -6240

$ws.Range("H126").Value = 4344.75
$ws.Range("I126").Value = 4344.75
$ws.Range("K126").Value = 13034.25
$ws.Range("M126").Value = -10564.25

$ws.Range("H132").Value = 10217140
$ws.Range("I132").Value = 9976.842000000001
$ws.Range("J132").Value = 34459150
$ws.Range("K132").Value = 29930.526
$ws.Range("L132").Value = 103377450
$ws.Range("M132").Value = -27400.526
$ws.Range("N132").Value = -103382510

$ws.Range("H138").Value = 169990
$ws.Range("J138").Value = 169990
$ws.Range("L138").Value = 169990
$ws.Range("N138").Value = -180270

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 11468.0625
$ws.Range("J7").Value = 13303.4
$ws.Range("L7").Value = 13303.4
$ws.Range("N7").Value = -13527.4

$ws.Range("H16").Value = 2069.7646
$ws.Range("I16").Value = 1829.4849
$ws.Range("K16").Value = 1829.4849
$ws.Range("M16").Value = -1659.4849

$ws.Range("I22").Value = 925.8
$ws.Range("J22").Value = 4761.5713
$ws.Range("K22").Value = 925.8
$ws.Range("L22").Value = 4761.5713
$ws.Range("M22").Value = -630.8
$ws.Range("N22").Value = -5351.5713

$ws.Range("I27").Value = 925.8
$ws.Range("J27").Value = 4761.5713
$ws.Range("K27").Value = 925.8
$ws.Range("L27").Value = 4761.5713
$ws.Range("M27").Value = -818.8
$ws.Range("N27").Value = -4975.5713

$ws.Range("H40").Value = 4255.6
$ws.Range("I40").Value = 4094.6667
$ws.Range("J40").Value = 4497
$ws.Range("K40").Value = 4094.6667
$ws.Range("L40").Value = 4497
$ws.Range("M40").Value = -3958.6667
$ws.Range("N40").Value = -4769

$ws.Range("H46").Value = 13374.777
$ws.Range("J46").Value = 1275
$ws.Range("L46").Value = 1275
$ws.Range("N46").Value = -1651

$ws.Range("H55").Value = 1000700.4
$ws.Range("I55").Value = 1786213.4
$ws.Range("K55").Value = 1786213.4
$ws.Range("M55").Value = -1786040.4

$ws.Range("H93").Value = 1817.625
$ws.Range("I93").Value = 1648.7142
$ws.Range("K93").Value = 1648.7142
$ws.Range("M93").Value = -400.7141999999999

$ws.Range("H125").Value = 88888
$ws.Range("J125").Value = 88888
$ws.Range("L125").Value = 88888
$ws.Range("N125").Value = -98728

$ws.Range("H126").Value = 11468.0625
$ws.Range("J126").Value = 13303.4
$ws.Range("L126").Value = 39910.2
$ws.Range("N126").Value = -44850.2

$ws.Range("H132").Value = 4696
$ws.Range("I132").Value = 4305.45
$ws.Range("J132").Value = 5296.846
$ws.Range("K132").Value = 12916.35
$ws.Range("L132").Value = 15890.538
$ws.Range("M132").Value = -10386.35
$ws.Range("N132").Value = -20950.538

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 41500
$ws.Range("J56").Value = 41500
$ws.Range("L56").Value = 41500
$ws.Range("N56").Value = -42928

$ws.Range("H62").Value = 3514.3635
$ws.Range("I62").Value = 3799.4
$ws.Range("K62").Value = 3799.4
$ws.Range("M62").Value = -3175.4

$ws.Range("H65").Value = 3514.3635
$ws.Range("I65").Value = 3799.4
$ws.Range("K65").Value = 18997
$ws.Range("M65").Value = -15877

$ws.Range("H74").Value = 12284.143
$ws.Range("J74").Value = 12452.728
$ws.Range("L74").Value = 12452.728
$ws.Range("N74").Value = -14324.728

$ws.Range("H75").Value = 39999.57
$ws.Range("I75").Value = 29999
$ws.Range("K75").Value = 29999
$ws.Range("M75").Value = -29063

$ws.Range("H77").Value = 12284.143
$ws.Range("J77").Value = 12452.728
$ws.Range("L77").Value = 37358.18399999999
$ws.Range("N77").Value = -46718.18399999999

$ws.Range("H78").Value = 39999.57
$ws.Range("I78").Value = 29999
$ws.Range("K78").Value = 89997
$ws.Range("M78").Value = -85317

$ws.Range("H82").Value = 41250
$ws.Range("I82").Value = 35000
$ws.Range("J82").Value = 60000
$ws.Range("K82").Value = 35000
$ws.Range("L82").Value = 60000
$ws.Range("M82").Value = -34617
$ws.Range("N82").Value = -60766

$ws.Range("H85").Value = 41250
$ws.Range("I85").Value = 35000
$ws.Range("J85").Value = 60000
$ws.Range("K85").Value = 35000
$ws.Range("L85").Value = 60000
$ws.Range("M85").Value = -33674
$ws.Range("N85").Value = -62652

$ws.Range("H107").Value = 2510.6428
$ws.Range("I107").Value = 1908.2858
$ws.Range("J107").Value = 3113
$ws.Range("K107").Value = 5724.857400000001
$ws.Range("L107").Value = 9339
$ws.Range("M107").Value = -3804.857400000001
$ws.Range("N107").Value = -13179

$ws.Range("H122").Value = 4426.364
$ws.Range("I122").Value = 3714
$ws.Range("K122").Value = 11142
$ws.Range("M122").Value = -8692

$ws.Range("H132").Value = 2370.524
$ws.Range("I132").Value = 1411.0834
$ws.Range("J132").Value = 3649.7778
$ws.Range("K132").Value = 4233.2502
$ws.Range("L132").Value = 10949.3334
$ws.Range("M132").Value = -1703.2502
$ws.Range("N132").Value = -16009.3334
